$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "position" column (E) for rows 4-7.
# Rows 4 and 5 (test2, test3) become "部门经理" (department manager).
$ws.Range("E4").Value = "部门经理"
$ws.Range("E5").Value = "部门经理"

# Rows 6 and 7 (test4, test5) become "员工" (employee).
$ws.Range("E6").Value = "员工"
$ws.Range("E7").Value = "员工"

# Move the active cell / selection to E11.
$null = $ws.Range("E11").Select()
